$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 474, shifting the existing rows 474:550 down
# to 475:551 (dimension grows from A1:R550 to A1:R551).
$ws.Rows.Item(474).Insert()

# Populate the newly inserted row with the new "Paine" / "1a (guarda)" record.
$ws.Cells.Item(474, 1).Value = 4
$ws.Cells.Item(474, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(474, 3).Value = "Los Lagos"
$ws.Cells.Item(474, 4).Value = 45218
$ws.Cells.Item(474, 5).Value = 10
$ws.Cells.Item(474, 6).Value = 100112045
$ws.Cells.Item(474, 7).Value = "Zapallo"
$ws.Cells.Item(474, 8).Value = "Paine"
$ws.Cells.Item(474, 9).Value = "1a (guarda)"
$ws.Cells.Item(474, 10).Value = 500
$ws.Cells.Item(474, 11).Value = 800
$ws.Cells.Item(474, 12).Value = 800
$ws.Cells.Item(474, 13).Value = 800
$ws.Cells.Item(474, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(474, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(474, 16).Value = 800
$ws.Cells.Item(474, 17).Value = 1
$ws.Cells.Item(474, 18).Value = "Hortaliza"
